# BOT; UPDATE DATA (#422)
# Append three new days of data (4/21, 4/22, 4/23 2020) to the "相談件数"
# sheet, correct the running totals on the last existing data row (4/20),
# move the footnote row down accordingly, and refresh the print area /
# selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")
$ws.Activate()

# --- Fix the running totals on row 86 (2020-04-20) ---------------------
$ws.Range("D86").Value = 229
$ws.Range("E86").Value = 5894

# --- Make room for 3 new data rows above the footnote row ---------------
# The footnote currently lives on row 87; inserting here shifts it (and
# copies row 86's formatting into the freshly inserted rows) down to 90.
$ws.Rows("87:89").Insert()

# --- Row 87: 2020-04-21 ---------------------------------------------------
$ws.Range("A87").Value = 43942
$ws.Range("B87").Value = 539
$ws.Range("C87").Value = 27289
$ws.Range("D87").Value = 174
$ws.Range("E87").Value = 6068

# --- Row 88: 2020-04-22 ---------------------------------------------------
$ws.Range("A88").Value = 43943
$ws.Range("B88").Value = 531
$ws.Range("C88").Value = 27820
$ws.Range("D88").Value = 127
$ws.Range("E88").Value = 6195

# --- Row 89: 2020-04-23 (counts not yet reported for B/D) -----------------
$ws.Range("A89").Value = 43944
$ws.Range("C89").Value = 27820
$ws.Range("E89").Value = 6195

# --- Update the print area defined name to cover the new rows -------------
$wb.Names.Item("Print_Area").RefersTo = '=相談件数!$A$1:$E$92'

# --- Move the view's active selection onto the newly added data -----------
$ws.Range("A89").Select()

$wb.Save()
